# Update column G (K = strikeouts) for rows 2-17 with recalculated values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 6
    4  = 1
    5  = 0
    6  = 0
    7  = 3
    8  = 1
    9  = 2
    10 = 0
    11 = 2
    12 = 2
    13 = 4
    14 = 3
    15 = 0
    16 = 5
    17 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
